# Apply ServerList.xlsx edits:
# - Update ServerIp values for existing rows (gateServer, logicServer)
# - Update NetThreadsNum values for existing rows
# - Add a new row for a "dbSercer" server entry

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 (gateServer, Id 10001): ServerIp 127.0.0.1 9998 -> 127.0.0.1 9997, NetThreadsNum 8 -> 4
$ws.Range("F5").Value = "127.0.0.1 9997"
$ws.Range("H5").Value = 4

# Row 6 (logicServer, Id 20001): ServerIp 127.0.0.1 9999 -> 127.0.0.1 9998, NetThreadsNum 4 -> 2
$ws.Range("F6").Value = "127.0.0.1 9998"
$ws.Range("H6").Value = 2

# Row 7 (new dbSercer entry, Id 30001)
$ws.Range("A7").Value = 30001
$ws.Range("B7").Value = 3
$ws.Range("C7").Value = "dbSercer"
$ws.Range("D7").Value = 1
$ws.Range("F7").Value = "127.0.0.1 9999"
$ws.Range("G7").Value = 2
$ws.Range("H7").Value = 2

# Update active selection to match the post-edit cursor position
[void]$ws.Range("O10").Select()
